$d = $word.ActiveDocument

# wdFindWrap / wdReplace constants used below:
#   Wrap   = 1  (wdFindContinue)
#   Replace = 2 (wdReplaceAll)

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# --- "En tant qu'utilisateur ... me connecter" -> "En tant que visiteur ... me connecter" ---
Replace-Text "En tant qu'utilisateur, je souhaite pouvoir me connecter en utilisant mon e-mail et mon mot de passe." "En tant que visiteur, je souhaite pouvoir me connecter en utilisant mon e-mail et mon mot de passe."

# --- "En tant qu'utilisateur ... réinitialiser" -> "En tant que visiteur ... réinitialiser" ---
Replace-Text "En tant qu’utilisateur, je souhaite pouvoir réinitialiser mon mot de passe en utilisant mon e-mail." "En tant que visiteur, je souhaite pouvoir réinitialiser mon mot de passe en utilisant mon e-mail."

# --- Merge the already-split "En tant qu" / "e visiteur" / ", je souhaite ... créer un compte ..." runs
#     back into a single run by re-writing the full sentence over the whole span. ---
Replace-Text "En tant que visiteur, je souhaite pouvoir créer un compte en utilisant mon e-mail depuis un lien sur la page de connexion." "En tant que visiteur, je souhaite pouvoir créer un compte en utilisant mon e-mail depuis un lien sur la page de connexion."

# --- "contactes" -> "contacts" fixes (French spelling correction) ---
Replace-Text "En tant qu’utilisateur, je souhaite gérer mes contactes et consulter les derniers messages reçus." "En tant qu’utilisateur, je souhaite gérer mes contacts et consulter les derniers messages reçus."

Replace-Text "pouvoir basculer sur la  discussion avec la le contacte sélectionné," "pouvoir basculer sur la  discussion avec la le contact sélectionné,"

Replace-Text "pouvoir visualiser le profil du contacte sélectionné," "pouvoir visualiser le profil du contact sélectionné,"

Replace-Text "pouvoir ajouter ou supprimer des contactes." "pouvoir ajouter ou supprimer des contacts."

Replace-Text "Gestion des contactes" "Gestion des contacts"

Replace-Text "En tant qu’utilisateur, je souhaite pouvoir gérer mes contactes." "En tant qu’utilisateur, je souhaite pouvoir gérer mes contacts."

Replace-Text "Étant donné que je suis sur la page de «Gestion des contactes»; je souhaite:" "Étant donné que je suis sur la page de «Gestion des contacts»; je souhaite:"

Replace-Text "ajouter un contacte depuis la liste des utilisateurs de l’application." "ajouter un contact depuis la liste des utilisateurs de l’application."

Replace-Text "supprimer un contacte de ma liste." "supprimer un contact de ma liste."

Replace-Text "Associer un contacte à cet événement." "Associer un contact à cet événement."

Replace-Text "associer une tache à un contacte." "associer une tache à un contact."

# --- styles.xml: Normal style paragraph format overflowPunct true -> false
#     (exposed in the Word object model as ParagraphFormat.HangingPunctuation) ---
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = $false

Write-Output "edit complete"
